$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 95, shifting existing rows 95-132 down to 96-133.
$ws.Rows.Item(95).Insert()

# Populate the newly inserted row 95 with the new data record.
$ws.Range("A95").Value = 1
$ws.Range("B95").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C95").Value = "Arica y Parinacota"
$ws.Range("D95").Value = 44455
$ws.Range("E95").Value = 15
$ws.Range("F95").Value = "Fruta"
$ws.Range("G95").Value = 100108
$ws.Range("H95").Value = "Tropicales y subtropicales"
$ws.Range("I95").Value = 100108006
$ws.Range("J95").Value = "Plátano"
$ws.Range("K95").Value = "Sin especificar"
$ws.Range("L95").Value = "Pintón"
$ws.Range("M95").Value = 120
$ws.Range("N95").Value = 19000
$ws.Range("O95").Value = 20000
$ws.Range("P95").Value = 19500
$ws.Range("Q95").Value = "$/caja 20 kilos"
$ws.Range("R95").Value = "Ecuador"
$ws.Range("S95").Value = 975
$ws.Range("T95").Value = 20
